{"js": "// Bold + yellow-highlight the \"Note that a 2-3 tree is a B-tree of order 3.\"\n// sentence in the B-tree section (the diff splits the run that used to\n// carry <w:br/> + the sentence into a plain <w:br/> run followed by a new\n// bold/highlighted run for the sentence text).\n\nconst body = context.document.body;\n\n// Locate the sentence (it's unique in the document).\nconst results = body.search(\"Note that a 2-3 tree is a B-tree of order 3.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target sentence not found\");\n}\n\nconst target = results.items[0];\ntarget.font.bold = true;\ntarget.font.highlightColor = \"yellow\";\nawait context.sync();\n\n// Word keeps a hidden \"_GoBack\" bookmark at the location of the most\n// recent edit. After editing the sentence above, it moves to collapse\n// right after the edited text, so re-create it there.\nconst goBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\ngoBack.load(\"isNullObject\");\nawait context.sync();\n\nif (!goBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n}\nconst endOfTarget = target.getRange(\"End\");\nendOfTarget.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Bold + yellow-highlight the \"Note that a 2-3 tree is a B-tree of order 3.\"\n# sentence in the B-tree section.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Note that a 2-3 tree is a B-tree of order 3.\")\n\nif ($found) {\n    # Use Range.Font (not Range directly) so the highlight only applies to\n    # the matched run instead of bleeding into the whole paragraph.\n    $rng.Font.Bold = 1\n    $rng.Font.HighlightColorIndex = 7   # wdYellow\n\n    # Word tracks the location of the most recent edit with a hidden\n    # \"_GoBack\" bookmark; after editing this sentence it collapses to\n    # right after the edited text, so recreate it there.\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $d.Bookmarks(\"_GoBack\").Delete()\n    }\n    $endRange = $d.Range($rng.End, $rng.End)\n    $d.Bookmarks.Add(\"_GoBack\", $endRange)\n}\n\nWrite-Output \"done\"\n"}
